# Auto-generated script applying scheduled market-data refresh to Zalera_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4723.613
$ws.Range("I11").Value = 4723.613
$ws.Range("K11").Value = 4723.613
$ws.Range("M11").Value = -4583.613
$ws.Range("H32").Value = 2300
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("M32").Value = -1674
$ws.Range("H100").Value = 29152.37
$ws.Range("I100").Value = 39383.77
$ws.Range("K100").Value = 39383.77
$ws.Range("M100").Value = -38842.77
$ws.Range("H116").Value = 6717
$ws.Range("I116").Value = 6499.75
$ws.Range("J116").Value = 6825.625
$ws.Range("K116").Value = 6499.75
$ws.Range("L116").Value = 6825.625
$ws.Range("M116").Value = -3057.75
$ws.Range("N116").Value = -13709.625
$ws.Range("H118").Value = 476.85715
$ws.Range("I118").Value = 473.33334
$ws.Range("K118").Value = 1420.00002
$ws.Range("M118").Value = 236.9999800000001
$ws.Range("H132").Value = 1408.3469
$ws.Range("I132").Value = 903.34283
$ws.Range("K132").Value = 2710.02849
$ws.Range("M132").Value = -180.0284900000001
$ws.Range("H137").Value = 125034710
$ws.Range("I137").Value = 250000000
$ws.Range("J137").Value = 69420
$ws.Range("K137").Value = 750000000
$ws.Range("L137").Value = 208260
$ws.Range("M137").Value = -749997450
$ws.Range("N137").Value = -213360
$ws.Range("H138").Value = 2973.195
$ws.Range("I138").Value = 1270.5
$ws.Range("J138").Value = 4944.737
$ws.Range("K138").Value = 3811.5
$ws.Range("L138").Value = 14834.211
$ws.Range("M138").Value = 1328.5
$ws.Range("N138").Value = -25114.211

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24928.234
$ws.Range("I32").Value = 24928.234
$ws.Range("K32").Value = 24928.234
$ws.Range("M32").Value = -24641.234
$ws.Range("H43").Value = 29749
$ws.Range("I43").Value = 20000
$ws.Range("J43").Value = 31698.8
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 31698.8
$ws.Range("M43").Value = -19687
$ws.Range("N43").Value = -32324.8
$ws.Range("H61").Value = 3066.5144
$ws.Range("I61").Value = 2614.9333
$ws.Range("K61").Value = 2614.9333
$ws.Range("M61").Value = -2402.9333
$ws.Range("H74").Value = 9246.120000000001
$ws.Range("I74").Value = 8587.85
$ws.Range("K74").Value = 8587.85
$ws.Range("M74").Value = -7713.85
$ws.Range("H77").Value = 9246.120000000001
$ws.Range("I77").Value = 8587.85
$ws.Range("K77").Value = 42939.25
$ws.Range("M77").Value = -38571.25
$ws.Range("H88").Value = 5396.2607
$ws.Range("J88").Value = 4856.1763
$ws.Range("L88").Value = 4856.1763
$ws.Range("N88").Value = -5668.1763
$ws.Range("H91").Value = 5396.2607
$ws.Range("J91").Value = 4856.1763
$ws.Range("L91").Value = 4856.1763
$ws.Range("N91").Value = -7664.1763
$ws.Range("H102").Value = 1873.0741
$ws.Range("I102").Value = 1925.9615
$ws.Range("J102").Value = 498
$ws.Range("K102").Value = 1925.9615
$ws.Range("L102").Value = 498
$ws.Range("M102").Value = -303.9614999999999
$ws.Range("N102").Value = -3742
$ws.Range("H110").Value = 22728800
$ws.Range("I110").Value = 31251488
$ws.Range("K110").Value = 31251488
$ws.Range("M110").Value = -31249443
$ws.Range("H136").Value = 3066.5144
$ws.Range("I136").Value = 2614.9333
$ws.Range("K136").Value = 7844.7999
$ws.Range("M136").Value = -5294.7999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 333417340
$ws.Range("I105").Value = 333417340
$ws.Range("K105").Value = 333417340
$ws.Range("M105").Value = -333415593
$ws.Range("H107").Value = 3564.16
$ws.Range("I107").Value = 1945.95
$ws.Range("K107").Value = 1945.95
$ws.Range("M107").Value = -25.95000000000005
$ws.Range("H134").Value = 4596.769
$ws.Range("I134").Value = 1711.1578
$ws.Range("J134").Value = 12429.143
$ws.Range("K134").Value = 5133.4734
$ws.Range("L134").Value = 37287.429
$ws.Range("M134").Value = -2598.4734
$ws.Range("N134").Value = -42357.429

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4234.381
$ws.Range("I16").Value = 2672.2856
$ws.Range("J16").Value = 5015.4287
$ws.Range("K16").Value = 2672.2856
$ws.Range("L16").Value = 5015.4287
$ws.Range("M16").Value = -2385.2856
$ws.Range("N16").Value = -5589.4287
$ws.Range("H22").Value = 669.9375
$ws.Range("I22").Value = 407.91666
$ws.Range("K22").Value = 407.91666
$ws.Range("M22").Value = -57.91665999999998
$ws.Range("H58").Value = 7534.4287
$ws.Range("I58").Value = 7197
$ws.Range("J58").Value = 7787.5
$ws.Range("K58").Value = 7197
$ws.Range("L58").Value = 7787.5
$ws.Range("M58").Value = -6994
$ws.Range("N58").Value = -8193.5
$ws.Range("H105").Value = 1087.1428
$ws.Range("I105").Value = 1087.1428
$ws.Range("K105").Value = 1087.1428
$ws.Range("M105").Value = 659.8571999999999
$ws.Range("H113").Value = 4234.381
$ws.Range("I113").Value = 2672.2856
$ws.Range("J113").Value = 5015.4287
$ws.Range("K113").Value = 2672.2856
$ws.Range("L113").Value = 5015.4287
$ws.Range("M113").Value = -502.2856000000002
$ws.Range("N113").Value = -9355.4287
$ws.Range("H134").Value = 5008.231
$ws.Range("I134").Value = 5029.1177
$ws.Range("J134").Value = 4866.2
$ws.Range("K134").Value = 15087.3531
$ws.Range("L134").Value = 14598.6
$ws.Range("M134").Value = -12552.3531
$ws.Range("N134").Value = -19668.6
$ws.Range("H136").Value = 7534.4287
$ws.Range("I136").Value = 7197
$ws.Range("J136").Value = 7787.5
$ws.Range("K136").Value = 21591
$ws.Range("L136").Value = 23362.5
$ws.Range("M136").Value = -19041
$ws.Range("N136").Value = -28462.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 11649.667
$ws.Range("I56").Value = 11649.667
$ws.Range("K56").Value = 11649.667
$ws.Range("M56").Value = -11119.667
$ws.Range("H57").Value = 2131
$ws.Range("I57").Value = 1446.5
$ws.Range("K57").Value = 4339.5
$ws.Range("M57").Value = -3780.5
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 2348.3333
$ws.Range("I41").Value = 439.7143
$ws.Range("K41").Value = 439.7143
$ws.Range("M41").Value = -84.71429999999998
$ws.Range("H70").Value = 12908.765
$ws.Range("I70").Value = 11266.385
$ws.Range("K70").Value = 11266.385
$ws.Range("M70").Value = -10996.385
$ws.Range("H73").Value = 12908.765
$ws.Range("I73").Value = 11266.385
$ws.Range("K73").Value = 11266.385
$ws.Range("M73").Value = -10330.385
$ws.Range("H80").Value = 2677.6667
$ws.Range("I80").Value = 2355.2
$ws.Range("J80").Value = 2908
$ws.Range("K80").Value = 2355.2
$ws.Range("L80").Value = 2908
$ws.Range("M80").Value = -1357.2
$ws.Range("N80").Value = -4904
$ws.Range("H83").Value = 2677.6667
$ws.Range("I83").Value = 2355.2
$ws.Range("J83").Value = 2908
$ws.Range("K83").Value = 11776
$ws.Range("L83").Value = 14540
$ws.Range("M83").Value = -6784
$ws.Range("N83").Value = -24524
$ws.Range("H132").Value = 3838.7334
$ws.Range("I132").Value = 1967.7693
$ws.Range("K132").Value = 5903.3079
$ws.Range("M132").Value = -3373.3079

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7150103
$ws.Range("I46").Value = 16669000
$ws.Range("K46").Value = 16669000
$ws.Range("M46").Value = -16668812
$ws.Range("H93").Value = 1866.4783
$ws.Range("J93").Value = 3456
$ws.Range("L93").Value = 3456
$ws.Range("N93").Value = -5952
$ws.Range("H100").Value = 8932799
$ws.Range("I100").Value = 12502869
$ws.Range("J100").Value = 7624.875
$ws.Range("K100").Value = 12502869
$ws.Range("L100").Value = 7624.875
$ws.Range("M100").Value = -12502328
$ws.Range("N100").Value = -8706.875
$ws.Range("H132").Value = 8052.4546
$ws.Range("I132").Value = 4882.7144
$ws.Range("K132").Value = 14648.1432
$ws.Range("M132").Value = -12118.1432
$ws.Range("H136").Value = 8854
$ws.Range("I136").Value = 7332
$ws.Range("K136").Value = 21996
$ws.Range("M136").Value = -19446

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H107").Value = 4928.4287
$ws.Range("I107").Value = 5500
$ws.Range("K107").Value = 16500
$ws.Range("M107").Value = -14580
$ws.Range("H132").Value = 6014.32
$ws.Range("I132").Value = 3241.182
$ws.Range("J132").Value = 8193.214
$ws.Range("K132").Value = 9723.545999999998
$ws.Range("L132").Value = 24579.642
$ws.Range("M132").Value = -7193.545999999998
$ws.Range("N132").Value = -29639.642
$ws.Range("H136").Value = 3323.76
$ws.Range("I136").Value = 2210.2632
$ws.Range("J136").Value = 6849.8335
$ws.Range("K136").Value = 6630.7896
$ws.Range("L136").Value = 20549.5005
$ws.Range("M136").Value = -4080.7896
$ws.Range("N136").Value = -25649.5005
